$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy the formatting of the
# existing header cell H1 (bold, bordered, centered/top-aligned) so the
# new headers look the same as B1:H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New data values for columns I (I0) and J (IF), rows 2-16.
$data = @(
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(3, 4),
    @(4, 5),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(3, 4),
    @(6, 7),
    @(9, 9),
    @(7, 7),
    @(5, 7),
    @(4, 5),
    @(3, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
